$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 165 (pushes the existing rows 165:197 down to 168:200),
# preserving formatting of the surrounding rows the same way Excel's
# Rows.Insert does.
$ws.Rows("165:167").Insert()

# Populate the 3 newly inserted rows with the new weekly price records.
# Columns A,B,C,E,F,G,H,I,J are constant for this market/product block.

# Row 165: Flame Seedless, Primera
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = 44995
$ws.Cells.Item(165, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100109
$ws.Cells.Item(165, 8).Value = "Uva"
$ws.Cells.Item(165, 9).Value = 100109001
$ws.Cells.Item(165, 10).Value = "Uva"
$ws.Cells.Item(165, 11).Value = "Flame Seedless"
$ws.Cells.Item(165, 12).Value = "Primera"
$ws.Cells.Item(165, 13).Value = 220
$ws.Cells.Item(165, 14).Value = 9000
$ws.Cells.Item(165, 15).Value = 10000
$ws.Cells.Item(165, 16).Value = 9545
$ws.Cells.Item(165, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(165, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(165, 19).Value = 530
$ws.Cells.Item(165, 20).Value = 18

# Row 166: Red Globe, Primera
$ws.Cells.Item(166, 1).Value = 11
$ws.Cells.Item(166, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(166, 3).Value = "Bíobío"
$ws.Cells.Item(166, 4).Value = 44995
$ws.Cells.Item(166, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(166, 5).Value = 8
$ws.Cells.Item(166, 6).Value = "Fruta"
$ws.Cells.Item(166, 7).Value = 100109
$ws.Cells.Item(166, 8).Value = "Uva"
$ws.Cells.Item(166, 9).Value = 100109001
$ws.Cells.Item(166, 10).Value = "Uva"
$ws.Cells.Item(166, 11).Value = "Red Globe"
$ws.Cells.Item(166, 12).Value = "Primera"
$ws.Cells.Item(166, 13).Value = 220
$ws.Cells.Item(166, 14).Value = 8000
$ws.Cells.Item(166, 15).Value = 9000
$ws.Cells.Item(166, 16).Value = 8545
$ws.Cells.Item(166, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(166, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(166, 19).Value = 475
$ws.Cells.Item(166, 20).Value = 18

# Row 167: Superior Seedless, Primera
$ws.Cells.Item(167, 1).Value = 11
$ws.Cells.Item(167, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(167, 3).Value = "Bíobío"
$ws.Cells.Item(167, 4).Value = 44995
$ws.Cells.Item(167, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(167, 5).Value = 8
$ws.Cells.Item(167, 6).Value = "Fruta"
$ws.Cells.Item(167, 7).Value = 100109
$ws.Cells.Item(167, 8).Value = "Uva"
$ws.Cells.Item(167, 9).Value = 100109001
$ws.Cells.Item(167, 10).Value = "Uva"
$ws.Cells.Item(167, 11).Value = "Superior Seedless"
$ws.Cells.Item(167, 12).Value = "Primera"
$ws.Cells.Item(167, 13).Value = 250
$ws.Cells.Item(167, 14).Value = 9000
$ws.Cells.Item(167, 15).Value = 10000
$ws.Cells.Item(167, 16).Value = 9600
$ws.Cells.Item(167, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(167, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(167, 19).Value = 533
$ws.Cells.Item(167, 20).Value = 18
